$wb = $excel.ActiveWorkbook

# --- Update the Date value on the Metadata sheet ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2025-07-21T11:52:46+00:00"

# --- Update the System URI values on the Include sheets ---
$wsInc0 = $wb.Worksheets.Item("Include #0")
$wsInc0.Range("B4").Value = "https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R20-Pays"

$wsInc1 = $wb.Worksheets.Item("Include #1")
$wsInc1.Range("B4").Value = "https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R268-PaysProvenanceISO"

$wsInc2 = $wb.Worksheets.Item("Include #2")
$wsInc2.Range("B4").Value = "https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R89-RegroupementPays"
